$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing designator "R2" (row 44), shifting rows below up.
$ws.Rows.Item(44).Delete()

# Update the named range to reflect the new, smaller extent of the data.
$wb.Names.Item("RP2040_Eins_top_pos").RefersTo = "=Sheet1!`$A`$1:`$E`$63"

# Update selection to match the recorded state after the edit.
$ws.Range("A44:XFD44").Select()

$wb.Save()
